# Rename sheet "Valve_30.0_600_3" -> "Valve_32.0_600_3" and make it the
# active sheet with cell D20 selected (moving the "active tab" / selection
# state off the previously-active sheet onto this one).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Valve_30.0_600_3")
$ws.Name = "Valve_32.0_600_3"

$ws.Activate()
$ws.Range("D20").Select()
